# Update countries & provincias Spain
#
# Refreshes the COVID "Pais" sheet with a newer data pull (22 May 2020,
# 10:35 instead of 10:05). Most rows only get new totals, but a handful of
# countries swapped rank relative to their neighbour since the last pull,
# so both the country name and its stats are rewritten for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Last updated" banner
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 10:35"

# Row 29: Singapur overtakes Portugal
$ws.Range("A29").Value = "Singapur"
$ws.Range("B29").Value = 30426
$ws.Range("C29").Value = 614
$ws.Range("D29").Value = 12117
$ws.Range("E29").Value = 18286
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 23

# Row 30: Portugal drops below Singapur
$ws.Range("A30").Value = "Portugal"
$ws.Range("B30").Value = 29912
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 6452
$ws.Range("E30").Value = 22183
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 1277

# Row 34: Polonia overtakes Indonesia & Ucrania
$ws.Range("A34").Value = "Polonia"
$ws.Range("B34").Value = 20379
$ws.Range("C34").Value = 236
$ws.Range("D34").Value = 8731
$ws.Range("E34").Value = 10675
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 973

# Row 35: Indonesia drops below Polonia
$ws.Range("A35").Value = "Indonesia"
$ws.Range("B35").Value = 20162
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 4838
$ws.Range("E35").Value = 14046
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 1278

# Row 36: Ucrania drops below Polonia
$ws.Range("A36").Value = "Ucrania"
$ws.Range("B36").Value = 20148
$ws.Range("C36").Value = 442
$ws.Range("D36").Value = 6585
$ws.Range("E36").Value = 12975
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 9
$ws.Range("H36").Value = 588

# Row 40: Rumania - stats refresh only (no reorder)
$ws.Range("A40").Value = "Rumania"
$ws.Range("B40").Value = 17585
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 10581
$ws.Range("E40").Value = 5845
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 3
$ws.Range("H40").Value = 1159

# Row 46: Filipinas - stats refresh only (no reorder)
$ws.Range("A46").Value = "Filipinas"
$ws.Range("B46").Value = 13597
$ws.Range("C46").Value = 163
$ws.Range("D46").Value = 3092
$ws.Range("E46").Value = 9648
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 11
$ws.Range("H46").Value = 857

# Row 52: Afganistan overtakes Chequia
$ws.Range("A52").Value = "Afganistan"
$ws.Range("B52").Value = 9216
$ws.Range("C52").Value = 540
$ws.Range("D52").Value = 996
$ws.Range("E52").Value = 8015
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 12
$ws.Range("H52").Value = 205

# Row 53: Chequia drops below Afganistan
$ws.Range("A53").Value = "Chequia"
$ws.Range("B53").Value = 8757
$ws.Range("C53").Value = 3
$ws.Range("D53").Value = 5932
$ws.Range("E53").Value = 2519
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 306

# Row 59: Malasia overtakes Australia
$ws.Range("A59").Value = "Malasia"
$ws.Range("B59").Value = 7137
$ws.Range("C59").Value = 78
$ws.Range("D59").Value = 5859
$ws.Range("E59").Value = 1163
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 115

# Row 60: Australia drops below Malasia
$ws.Range("A60").Value = "Australia"
$ws.Range("B60").Value = 7095
$ws.Range("C60").Value = 14
$ws.Range("D60").Value = 6478
$ws.Range("E60").Value = 516
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 101

# Row 65: Oman - stats refresh only (no reorder)
$ws.Range("A65").Value = "Oman"
$ws.Range("B65").Value = 6370
$ws.Range("C65").Value = 0
$ws.Range("D65").Value = 1821
$ws.Range("E65").Value = 4517
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 1
$ws.Range("H65").Value = 32

# Row 93: Lituania overtakes Somalia
$ws.Range("A93").Value = "Lituania"
$ws.Range("B93").Value = 1604
$ws.Range("C93").Value = 11
$ws.Range("D93").Value = 1111
$ws.Range("E93").Value = 432
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 61

# Row 94: Somalia drops below Lituania
$ws.Range("A94").Value = "Somalia"
$ws.Range("B94").Value = 1594
$ws.Range("C94").Value = 0
$ws.Range("D94").Value = 204
$ws.Range("E94").Value = 1329
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 61
